$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New matchup rows (fall 13 week 3 inputs) to append starting at row 2226
$data = @(
    @(3,0,3,3),
    @(5,2,5,0),
    @(6,2,5,0),
    @(5,2,4,0),
    @(3,2,4,0),
    @(4,0,6,2),
    @(3,1,3,2),
    @(4,0,4,2),
    @(5,3,7,0),
    @(4,3,3,0),
    @(4,2,4,0),
    @(3,0,3,3),
    @(5,0,6,2),
    @(3,0,3,3),
    @(5,2,7,0),
    @(5,1,4,2),
    @(4,0,7,3),
    @(6,0,5,2),
    @(4,2,2,1),
    @(4,0,4,2),
    @(2,1,2,2),
    @(6,2,5,1),
    @(6,1,7,2),
    @(3,3,3,0)
)

$startRow = 2226
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$lastRow = $startRow + $data.Count - 1
$lastCell = "A" + ($lastRow + 1)
$ws.Range($lastCell).Select() | Out-Null
